# Conversão para TestNG e adição do report
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the username value in A2 (shared string "joao1211" -> "joao12366")
$ws.Range("A2").Value2 = "joao12366"

# Header row (A1:L1): center the text both horizontally and vertically
$headerRng = $ws.Range("A1:L1")
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4108

# Row 6 (A6): add a thin right-edge divider so it matches the header box style
$ws.Range("A6").Borders.Item(10).LineStyle = 1
$ws.Range("A6").Borders.Item(10).Weight = 2

# Move the active selection to A3
$ws.Range("A3").Select()
